$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDN")
$ws.Columns("D:D").Insert()
